# Re-run of the data preprocessing / cleaning script ("Preprocesamiento y
# Limpieza de Datos"). The upstream source was re-fetched, so a handful of
# rows picked up refreshed population/area/capital/subregion values, the
# dependent population_density / language_density columns were
# recomputed, and every row's processing timestamp was refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Refreshed source values (population, area, capital, subregion) for
#    the rows where the re-fetched data differs from the previous run.
# ---------------------------------------------------------------------

# Row 81 - Mayotte: capital resolved
$ws.Range("J81").Value = '["Mamoudzou"]'

# Row 92 - El Salvador: population refreshed
$ws.Range("G92").Value = 6486201

# Row 110 - Canada: subregion resolved
$ws.Range("F110").Value = "North America"

# Row 115 - Spain: capital resolved
$ws.Range("J115").Value = '["Madrid"]'

# Row 116 - Slovenia: population refreshed
$ws.Range("G116").Value = 2099235

# Row 118 - Saint Pierre and Miquelon: capital lookup failed this run
$ws.Range("J118").Value = "Unknown"

# Row 148 - Cameroon: area refreshed
$ws.Range("H148").Value = 475442

# Row 151 - Suriname: population refreshed
$ws.Range("G151").Value = 586634

# Row 165 - Romania: subregion lookup failed this run
$ws.Range("F165").Value = "Unknown"

# Row 171 - Vatican City: subregion lookup failed this run
$ws.Range("F171").Value = "Unknown"

# Row 182 - Norway: subregion resolved
$ws.Range("F182").Value = "Northern Europe"

# Row 188 - South Korea: capital lookup failed this run
$ws.Range("J188").Value = "Unknown"

# Row 192 - United States Minor Outlying Islands: area refreshed
$ws.Range("H192").Value = 34.80924866606932

# Row 196 - Moldova: area refreshed
$ws.Range("H196").Value = 33846

# Row 242 - Saint Martin: population refreshed
$ws.Range("G242").Value = 39039

# Row 249 - New Zealand: area refreshed
$ws.Range("H249").Value = 268208.9349268354

# ---------------------------------------------------------------------
# 2) Recompute the dependent columns for every row whose population (G)
#    or area (H) changed:
#       O = population_density   = population / area
#       U = language_density     = language_count / population * 1e6
# ---------------------------------------------------------------------

$recalcRows = @(92, 116, 148, 151, 192, 196, 242, 249)
foreach ($r in $recalcRows) {
    $population = $ws.Cells.Item($r, 7).Value   # column G
    $area       = $ws.Cells.Item($r, 8).Value   # column H
    $langCount  = $ws.Cells.Item($r, 19).Value  # column S

    $ws.Cells.Item($r, 15).Value = $population / $area
    $ws.Cells.Item($r, 21).Value = $langCount / $population * 1000000
}

# ---------------------------------------------------------------------
# 3) Refresh the per-row processing timestamp (column N) to reflect the
#    time this run of the script executed.
# ---------------------------------------------------------------------

$ws.Range("ZZ1").Formula = '=TEXT(NOW(),"yyyy-mm-dd")'
$datePart = $ws.Range("ZZ1").Text
$ws.Range("ZZ1").Formula = '=TEXT(NOW(),"hh:mm:ss")'
$timePart = $ws.Range("ZZ1").Text
$ws.Range("ZZ1").ClearContents()

$stamp = $datePart + "T" + $timePart + ".000000"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 14).Value = $stamp
}
